$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for the numeric-looking "Price" column (D) so values
# are stored as text (matching the original inline-string cells), not numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.827.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.346.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "471.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.502"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.346.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0958"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.318"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.757.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "54.883.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.353.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "311.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "55.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.390"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.449.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0747"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.816"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0945"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0522"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "251.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.775.40"
$ws.Range("D51").Style = "Normal"

# Remaining text / volume-percentage columns are not numeric-parseable,
# so plain assignment keeps them as text.
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("E3").Value = "  -5.39%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  -3.65%  "
$ws.Range("E6").Value = "  -2.55%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("E9").Value = "  -5.98%  "
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("E11").Value = "  -6.45%  "
$ws.Range("E12").Value = "  -4.59%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("E14").Value = "  -5.38%  "
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("E16").Value = "  -5.98%  "
$ws.Range("E17").Value = "  -5.10%  "
$ws.Range("E18").Value = "  -5.81%  "
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("E21").Value = "  -5.45%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  -4.17%  "
$ws.Range("E24").Value = "  -4.91%  "
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("E27").Value = "  -5.96%  "
$ws.Range("E28").Value = "  -5.52%  "
$ws.Range("E29").Value = "  -6.55%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -5.60%  "
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("E35").Value = "  -4.00%  "
$ws.Range("E36").Value = "  -5.36%  "
$ws.Range("E37").Value = "  -4.87%  "
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  -4.76%  "
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("E44").Value = "  -5.85%  "
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("E49").Value = "  -8.99%  "
$ws.Range("E50").Value = "  -5.39%  "
$ws.Range("E51").Value = "  -6.07%  "
